$d = $word.ActiveDocument

# Update the date heading in the first paragraph
$heading = $d.Paragraphs.Item(1).Range
$heading.Text = "2025-06-22 Sunday"

# Update each cell of the answer table (20 rows x 5 cols) positionally,
# since several old values repeat and a global Find/Replace would be ambiguous.
$t = $d.Tables.Item(1)
$answers = @(
    @("12+59=71", "3+38=41", "68+7=75", "35+29=64", "9+39=48"),
    @("16-7=9", "84-37=47", "17+29=46", "55+29=84", "91-88=3"),
    @("33-29=4", "66+27=93", "23+69=92", "32-16=16", "61-27=34"),
    @("56+5=61", "84-55=29", "24-7=17", "48+39=87", "81-65=16"),
    @("61-44=17", "69+24=93", "44+8=52", "35-9=26", "94-39=55"),
    @("29+4=33", "91-78=13", "80-42=38", "90-61=29", "66-49=17"),
    @("64-57=7", "97-18=79", "47+5=52", "15+58=73", "74-49=25"),
    @("66-48=18", "85-66=19", "59+15=74", "69+4=73", "25+57=82"),
    @("50-16=34", "92-9=83", "11-2=9", "69+6=75", "60-13=47"),
    @("9+83=92", "54+9=63", "32-4=28", "91-66=25", "19+7=26"),
    @("29+54=83", "3+78=81", "91-19=72", "85-57=28", "81-34=47"),
    @("15+48=63", "98-19=79", "59+18=77", "29+38=67", "55-7=48"),
    @("12+19=31", "69+14=83", "75+9=84", "74-59=15", "35+27=62"),
    @("60-36=24", "73-6=67", "81-18=63", "25+46=71", "70-16=54"),
    @("85-18=67", "98-9=89", "40-39=1", "77+15=92", "5+89=94"),
    @("38-19=19", "63-29=34", "71-9=62", "45-7=38", "81-29=52"),
    @("48+28=76", "9+26=35", "8+43=51", "51-16=35", "54-39=15"),
    @("23+28=51", "34+28=62", "8+19=27", "55+16=71", "66+16=82"),
    @("72-27=45", "13+59=72", "77-28=49", "78+9=87", "20-12=8"),
    @("72-64=8", "70-51=19", "9+39=48", "46+35=81", "55-39=16")
)

for ($r = 1; $r -le $t.Rows.Count; $r++) {
    $rowVals = $answers[$r - 1]
    for ($c = 1; $c -le $t.Rows.Item($r).Cells.Count; $c++) {
        $cell = $t.Rows.Item($r).Cells.Item($c)
        $cell.Range.Text = $rowVals[$c - 1]
    }
}

Write-Output "done"